$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (header "K") values regenerated (strikeouts recalculated per
# commit message "regen save_data to use K instead of Strike#, regen
# std/mean, calc and write s_vals"). Rows 32 and 42 are unchanged.
$newK = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 2
    6  = 2
    7  = 0
    8  = 0
    9  = 0
    10 = 1
    11 = 1
    12 = 1
    13 = 0
    14 = 0
    15 = 1
    16 = 1
    17 = 2
    18 = 2
    19 = 1
    20 = 0
    21 = 1
    22 = 1
    23 = 1
    24 = 0
    25 = 1
    26 = 3
    27 = 2
    28 = 0
    29 = 3
    30 = 1
    31 = 1
    33 = 0
    34 = 1
    35 = 2
    36 = 2
    37 = 3
    38 = 1
    39 = 0
    40 = 1
    41 = 1
    43 = 1
    44 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
